$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.421.06'
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("D3").Value = '1.665.80'
$ws.Range("E3").Value = '  +1.21%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("D5").Value = '312.78'
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = '0.3961'
$ws.Range("E7").Value = '  +1.74%  '
$ws.Range("D8").Value = '0.3938'
$ws.Range("D9").Value = '52.17'
$ws.Range("E9").Value = '  +6.68%  '
$ws.Range("D10").Value = '1.395'
$ws.Range("E10").Value = '  +3.69%  '
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '0.08574'
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").Value = '24.51'
$ws.Range("E13").Value = '  +3.19%  '
$ws.Range("D14").Value = '7.324'
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("D15").Value = '7.934'
$ws.Range("E15").Value = '  +6.17%  '
$ws.Range("E16").Value = '  +4.56%  '
$ws.Range("D17").Value = '1.666.99'
$ws.Range("E17").Value = '  +1.76%  '
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").Value = '0.06991'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '20.61'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '6.998'
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("D22").Value = '0.9981'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = '13.76'
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("D24").Value = '24.426.60'
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").Value = '3.121'
$ws.Range("E25").Value = '  +15.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.430'
$ws.Range("E26").Value = '  +4.40%  '
$ws.Range("D27").Value = '22.54'
$ws.Range("D28").Value = '158.07'
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").Value = '142.81'
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("D30").Value = '5.424'
$ws.Range("E30").Value = '  +3.34%  '
$ws.Range("D31").Value = '8.064'
$ws.Range("E31").Value = '  -6.95%  '
$ws.Range("D32").Value = '2.533'
$ws.Range("E32").Value = '  +2.96%  '
$ws.Range("D33").Value = '1.851.48'
$ws.Range("E33").Value = '  +1.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.070'
$ws.Range("E34").Value = '  +11.65%  '
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = '0.03064'
$ws.Range("E35").Value = '  +5.51%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.08266'
$ws.Range("E36").Value = '  +3.19%  '
$ws.Range("D37").Value = '6.925'
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("D38").Value = '11.14'
$ws.Range("E38").Value = '  +12.32%  '
$ws.Range("D39").Value = '0.2765'
$ws.Range("E39").Value = '  +2.85%  '
$ws.Range("D40").Value = '0.09251'
$ws.Range("E40").Value = '  +0.60%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.7717'
$ws.Range("E41").Value = '  +1.72%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '13.82'
$ws.Range("E42").Value = '  +5.89%  '
$ws.Range("D43").Value = '1.451'
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").Value = '16.57'
$ws.Range("E44").Value = '  +3.70%  '
$ws.Range("D45").Value = '0.7127'
$ws.Range("E45").Value = '  +3.48%  '
$ws.Range("D46").Value = '2.544'
$ws.Range("E46").Value = '  +2.83%  '
$ws.Range("D47").Value = '4.138'
$ws.Range("E47").Value = '  +1.10%  '
$ws.Range("D48").Value = '0.9999'
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("D49").Value = '0.08432'
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("D50").Value = '136.81'
$ws.Range("E50").Value = '  +2.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.270'
$ws.Range("E51").Value = '  +1.07%  '
